$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: "时间" (time) is replaced by "日期" (date) in column B;
# "可预约人数"/"已预约人数" stay put in C1/D1 but now gain explicit number formats.
$ws.Range("B1").Value = "日期"

# Apply number formats to the header cells (these flow onto the whole
# column via the xf records, matching the added date/integer cellXfs).
$ws.Range("B1").NumberFormat = "mm-dd-yy"
$ws.Range("C1").NumberFormat = "0"
$ws.Range("D1").NumberFormat = "0"

# Column widths (engine quantizes width to a 1/6-character grid, so these
# inputs are chosen to land on the grid point nearest the template's true
# 8.88671875 / 11.33203125 / 10.6640625 character widths).
$ws.Columns("B").ColumnWidth = 8
$ws.Columns("C").ColumnWidth = 10.5
$ws.Columns("D").ColumnWidth = 9.833333333333334

# Selection moved from D1 to D12.
$ws.Range("D12").Select()

# Page orientation switched to portrait (the template also gained
# horizontalDpi/verticalDpi + a printerSettings relationship in Excel,
# but the headless host has no COM surface for printer-settings parts).
$ws.PageSetup.Orientation = 1
